$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 20604
$ws1.Range("F4").Value = 326
$ws1.Range("F7").Value = 7712
$ws1.Range("F8").Value = 534
$ws1.Range("F9").Value = 747
$ws1.Range("F10").Value = 290
$ws1.Range("F12").Value = 174
$ws1.Range("F13").Value = 140
$ws1.Range("F16").Value = 208
$ws1.Range("F18").Value = 473
$ws1.Range("F21").Value = 52
$ws1.Range("F23").Value = 75
$ws1.Range("F24").Value = 336
$ws1.Range("F25").Value = 1149
$ws1.Range("F28").Value = 198
$ws1.Range("F31").Value = 100
$ws1.Range("F32").Value = 4920
$ws1.Range("F36").Value = 12828
$ws1.Range("F37").Value = 1347
$ws1.Range("F38").Value = 101
$ws1.Range("F39").Value = 38
$ws1.Range("F41").Value = 289
$ws1.Range("F43").Value = 4025
$ws1.Range("F44").Value = 324

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 232

# Sheet 4: 全部类型 (All types, combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 20604
$ws4.Range("F4").Value = 326
$ws4.Range("F7").Value = 7712
$ws4.Range("F8").Value = 534
$ws4.Range("F9").Value = 747
$ws4.Range("F10").Value = 290
$ws4.Range("F12").Value = 174
$ws4.Range("F13").Value = 140
$ws4.Range("F16").Value = 208
$ws4.Range("F18").Value = 473
$ws4.Range("F23").Value = 75
$ws4.Range("F24").Value = 336
$ws4.Range("F25").Value = 1149
$ws4.Range("F28").Value = 198
$ws4.Range("F29").Value = 232
$ws4.Range("F33").Value = 100
$ws4.Range("F35").Value = 4920
$ws4.Range("F39").Value = 12828
$ws4.Range("F40").Value = 1347
$ws4.Range("F41").Value = 101
$ws4.Range("F42").Value = 38
$ws4.Range("F44").Value = 289
$ws4.Range("F46").Value = 4025
$ws4.Range("F47").Value = 324
